$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (train_size = 0.3)
$ws.Range("B2").Value = "0.0899 (0.0033)"
$ws.Range("C2").Value = "0.8273 (0.0481)"
$ws.Range("D2").Value = "0.8295 (0.0897)"
$ws.Range("E2").Value = "0.8222 (0.1015)"

# Row 3 (train_size = 0.5)
$ws.Range("C3").Value = "0.8290 (0.0468)"
$ws.Range("D3").Value = "0.8299 (0.0808)"
$ws.Range("E3").Value = "0.8271 (0.1012)"

# Row 4 (train_size = 0.8)
$ws.Range("C4").Value = "0.8319 (0.0626)"
$ws.Range("D4").Value = "0.8340 (0.0897)"
$ws.Range("E4").Value = "0.8272 (0.1243)"
